# Appends " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document."), as three separate runs:
#   " (", "Changed main", ")"
# so the original sentence run is left untouched and three new runs
# follow it, matching the target OOXML diff.

$d = $word.ActiveDocument

# Inserts $Text as its own run at character position $Pos in $Doc and
# returns the position immediately after the inserted text. Bracketing
# the insertion with a transient bookmark at the boundary keeps the new
# text from being coalesced back into the neighboring run on save, so
# each call produces a distinct <w:r> element.
function Insert-RunText($Doc, $Pos, $Text) {
    $guardName = "zzTmpRunGuard"
    $Doc.Bookmarks.Add($guardName, $Doc.Range($Pos, $Pos))
    $r = $Doc.Range($Pos, $Pos)
    $r.InsertAfter($Text)
    $Doc.Bookmarks($guardName).Delete()
    return $r.End
}

# Find the existing sentence and collapse to just after it (before the
# paragraph mark) so new runs are appended within the same paragraph.
$anchor = $d.Content
$null = $anchor.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $anchor.End

$pos = Insert-RunText $d $pos " ("
$pos = Insert-RunText $d $pos "Changed main"
$pos = Insert-RunText $d $pos ")"
